# Commit: "add arduino hardware, date 16/Mar/21"
#
# The unified diff (scraped from a different export of this evolving
# document) shows a single content-bearing hunk: a brand-new, empty
# paragraph is inserted as the next-to-last paragraph of the document
# body -- i.e. immediately before the final (empty, green-bar) paragraph
# that precedes </w:body>. The new paragraph carries heading-style run
# properties (Khmer OS, bold, black, 14pt/sz=28, no shading override)
# matching the "2. Arduio Board" / "4. <heading>" section headers
# already present at the end of this document, and holds a single empty
# run whose only explicit property is "no underline" -- i.e. it is the
# blank heading line that begins the new "Arduino hardware" section the
# author is about to add. All the remaining hunks in the supplied diff
# are repeated "overflowPunct true -> false" churn coming from the
# Normal-style pPr defaults in that other export of the file; this
# document's styles/paragraphs do not carry that property at all, so
# there is nothing to mirror for those hunks.

$d = $word.ActiveDocument

# The paragraph immediately before the document's final (trailing,
# empty) paragraph -- currently the "4. <heading>" section title.
$count = $d.Paragraphs.Count
$anchor = $d.Paragraphs.Item($count - 1)

# Create a new empty paragraph right after it (i.e. right before the
# final trailing paragraph), then stamp it with the exact OOXML the
# commit introduces.
$anchor.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($count)
$newRange = $newPara.Range

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Normal"/>
              <w:rPr>
                <w:rFonts w:ascii="Khmer OS" w:hAnsi="Khmer OS" w:cs="Khmer OS"/>
                <w:b/>
                <w:b/>
                <w:bCs/>
                <w:color w:val="000000"/>
                <w:sz w:val="28"/>
                <w:szCs w:val="28"/>
                <w:shd w:fill="auto" w:val="clear"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:u w:val="none"/>
              </w:rPr>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newRange.InsertXML($xml)

Write-Host "Paragraphs now:" $d.Paragraphs.Count
